$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old columns C and D entirely (they are being removed)
$ws.Range("C1:D2").Clear()

# Update headers
$ws.Range("A1").Value = "Position"
$ws.Range("B1").Value = "Value"

# Update data rows
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0
